$wb = $excel.ActiveWorkbook

# The "Repayment schedule" sheet gets a new column inserted before column N
# (i.e. a new, blank "Late" spacer column), and becomes the active/selected sheet
# (previously "NewLoanInput" was the selected tab).
$wsRepay = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (shifts N:P -> O:Q), carrying
# over the column width from the column immediately to the left (M), as
# Excel does for a plain column insert.
$leftWidth = $wsRepay.Columns("M").ColumnWidth
$wsRepay.Columns("N").Insert()
$wsRepay.Columns("N").ColumnWidth = $leftWidth

# Update the selection on the Repayment schedule sheet
$wsRepay.Range("R7").Select()

# Make "Repayment schedule" the active sheet/tab (tabSelected moves here)
$wsRepay.Activate()
$wsRepay.Select()
